# feat: add 2022-Q3 data
#
# 1) Insert a new "2022-Q3" quarter worksheet (copied from the "2022-Q2"
#    sheet so the header/index-column formatting matches), positioned right
#    before "2022-Q2", and fill it with the new quarter's fund data.
# 2) Insert a new row at the top of the "总计" (summary) sheet's data with
#    the 2022-Q3 totals; every other row shifts down by one, values unchanged.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the "2022-Q2" sheet to create "2022-Q3" right before it
# ---------------------------------------------------------------------
$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($srcQ2)
$ws2 = $wb.Worksheets.Item("2022-Q2 (2)")
$ws2.Name = "2022-Q3"

# The copied sheet only has 10 data rows (rows 2-11); the new quarter needs
# an 11th data row (row 12), so extend formatting down one more row.
$ws2.Range("A11:H11").Copy()
$ws2.Range("A12:H12").PasteSpecial(-4122)

# 2022-Q3 fund holdings: fund code, fund name, fund scale, total equity
# position, position ratio, held market value, position rank.
$fundData = @(
  @("161222", "国投瑞银瑞利灵活配置混合（LOF）A", "23.90", "68.57", "2.38", "0.5688", 5),
  @("010338", "国投瑞银远见成长混合A",             "10.51", "88.69", "2.48", "0.2606", 9),
  @("121010", "国投瑞银瑞源灵活配置混合A",         "9.45",  "70.88", "2.36", "0.2230", 5),
  @("015652", "国投瑞银瑞利灵活配置混合（LOF）C", "8.00",  "68.57", "2.38", "0.1904", 5),
  @("005904", "华泰保兴成长优选混合A",             "4.00",  "80.96", "3.74", "0.1496", 10),
  @("012132", "华泰保兴价值成长混合A",             "0.95",  "85.32", "4.81", "0.0457", 10),
  @("010339", "国投瑞银远见成长混合C",             "1.75",  "88.69", "2.48", "0.0434", 9),
  @("015572", "国投瑞银瑞源灵活配置混合C",         "1.29",  "70.88", "2.36", "0.0304", 5),
  @("620001", "金元顺安宝石动力混合",               "1.01",  "40.12", "2.88", "0.0291", 6),
  @("012177", "华泰保兴价值成长混合C",             "0.11",  "85.32", "4.81", "0.0053", 10),
  @("005905", "华泰保兴成长优选混合C",             "0.14",  "80.96", "3.74", "0.0052", 10)
)

# Columns B, D, E, F, G are numeric-looking text in the source data (e.g.
# fund code "005904", fund scale "23.90") so they are entered with a
# leading "'" to force text storage (preserving leading/trailing zeros),
# then QuotePrefix is cleared so the saved cell carries no extra style -
# same shape as the source sheets (plain inline text, no quote marker).
$r = 2
foreach ($row in $fundData) {
  $ws2.Cells.Item($r, 1).Value = $r - 2

  $ws2.Cells.Item($r, 2).Value = "'" + $row[0]
  $ws2.Cells.Item($r, 2).QuotePrefix = $false

  $ws2.Cells.Item($r, 3).Value = $row[1]

  $ws2.Cells.Item($r, 4).Value = "'" + $row[2]
  $ws2.Cells.Item($r, 4).QuotePrefix = $false

  $ws2.Cells.Item($r, 5).Value = "'" + $row[3]
  $ws2.Cells.Item($r, 5).QuotePrefix = $false

  $ws2.Cells.Item($r, 6).Value = "'" + $row[4]
  $ws2.Cells.Item($r, 6).QuotePrefix = $false

  $ws2.Cells.Item($r, 7).Value = "'" + $row[5]
  $ws2.Cells.Item($r, 7).QuotePrefix = $false

  $ws2.Cells.Item($r, 8).Value = $row[6]

  $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: insert the new 2022-Q3 summary row into "总计", pushing the
# existing rows down, then rewrite the table so the running index in
# column A stays sequential (0,1,2,...) and every quarter's stats land on
# the right row.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

# Extend formatting (bold/bordered index cell in column A, plain cells for
# B:D) down to the newly-needed last row (row 9) by copying row 8's format.
$wsTotal.Range("A8:D8").Copy()
$wsTotal.Range("A9:D9").PasteSpecial(-4122)

# The freshly-inserted row 2 picked up the bold header formatting from row
# 1; re-format it like the other data rows (row 3) instead.
$wsTotal.Range("A3:D3").Copy()
$wsTotal.Range("A2:D2").PasteSpecial(-4122)

# Rewrite the whole data block (rows 2-9) with final values, newest quarter
# first, so column A is the simple sequential index 0..7.
$summaryData = @(
  @("2022-Q3", 11, 1.55),
  @("2022-Q2", 10, 2.19),
  @("2022-Q1", 19, 3.77),
  @("2021-Q4", 17, 4.26),
  @("2021-Q3", 7, 0.86),
  @("2021-Q2", 5, 0.31),
  @("2021-Q1", 8, 0.68),
  @("2020-Q4", 6, 0.52)
)

$r = 2
foreach ($row in $summaryData) {
  $wsTotal.Cells.Item($r, 1).Value = $r - 2
  $wsTotal.Cells.Item($r, 2).Value = $row[0]
  $wsTotal.Cells.Item($r, 3).Value = $row[1]
  $wsTotal.Cells.Item($r, 4).Value = $row[2]
  $r = $r + 1
}

Write-Host "2022-Q3 sheet inserted and summary sheet updated"
